$d = $word.ActiveDocument

# --- Edit 1: "autre," -> "aultre," ---
$d.Content.Find.Execute("autre,", $true, $false, $false, $false, $false, $true, 1, $false, "aultre,", 2) | Out-Null

# --- Edit 2: move the "a" run from after "<corr>" to before it ---
# i.e. "<corr>affin</corr>" becomes "a<corr>ffin</corr>"
# Anchor the search right after the text we just fixed, to make sure we
# land on the correct occurrence of "<corr>affin</corr>" in the document.
$anchor = $d.Content
$anchor.Find.Execute("aultre,", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null

$target = $d.Range($anchor.End, $d.Content.End)
$target.Find.Execute("<corr>affin</corr>", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$corrStart = $target.Start

# The run containing just "a" sits right after the 6 characters of "<corr>"
$aRun = $d.Range($corrStart + 6, $corrStart + 7)

# Copy it (with its run formatting) to a new location right before "<corr>"
$dest = $d.Range($corrStart, $corrStart)
$dest.FormattedText = $aRun.FormattedText

# Remove the original "a" (now shifted one character to the right because
# of the text we just inserted)
$aRunOld = $d.Range($corrStart + 1 + 6, $corrStart + 1 + 7)
$aRunOld.Delete() | Out-Null
